$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Cell B11 currently holds the text "R40"; change it to the text "1".
# Assigning a plain numeric-looking string (e.g. "1") would make Excel
# store it as a number, so it is first written as a text formula and
# then converted to a static value, which keeps the cell's existing
# style and "stored as text" type intact.
$cell = $ws.Range("B11")
$cell.Formula = '="1"'
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null  # xlPasteValues
$excel.CutCopyMode = $false
